$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '22.446.53'
$ws.Range("E2").Value = '  +9.02%  '
$ws.Range("D3").Value = '1.607.32'
$ws.Range("E3").Value = '  +8.78%  '
$ws.Range("E4").Value = '  -0.89%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.9912'
$ws.Range("E5").Value = '  +3.47%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '303.57'
$ws.Range("E6").Value = '  +8.24%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3679'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3408'
$ws.Range("E8").Value = '  +11.05%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '42.29'
$ws.Range("E9").Value = '  +5.58%  '
$ws.Range("E10").Value = '  +7.23%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07086'
$ws.Range("E11").Value = '  +6.09%  '
$ws.Range("E12").Value = '  -0.62%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '19.75'
$ws.Range("E13").Value = '  +9.29%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.931'
$ws.Range("E14").Value = '  +7.31%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.658'
$ws.Range("E15").Value = '  +6.98%  '
$ws.Range("E16").Value = '  +5.28%  '
$ws.Range("D17").Value = '1.602.80'
$ws.Range("E17").Value = '  +8.46%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.9916'
$ws.Range("E18").Value = '  +3.39%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06837'
$ws.Range("E19").Value = '  +14.89%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '78.16'
$ws.Range("E20").Value = '  +11.66%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.057'
$ws.Range("E21").Value = '  +9.99%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '16.13'
$ws.Range("E22").Value = '  +11.60%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '11.87'
$ws.Range("E23").Value = '  +7.24%  '
$ws.Range("D24").Value = '22.468.76'
$ws.Range("E24").Value = '  +8.91%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.379'
$ws.Range("E25").Value = '  +5.25%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.546'
$ws.Range("E26").Value = '  +20.60%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '151.20'
$ws.Range("E27").Value = '  +5.55%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '19.61'
$ws.Range("E28").Value = '  +13.39%  '
$ws.Range("D29").Value = '1.784.32'
$ws.Range("E29").Value = '  +8.87%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '121.28'
$ws.Range("E30").Value = '  +6.46%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.103'
$ws.Range("E31").Value = '  +3.24%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.093'
$ws.Range("E32").Value = '  +21.12%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.9519'
$ws.Range("E33").Value = '  +17.00%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.08287'
$ws.Range("E34").Value = '  +3.99%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.641'
$ws.Range("E35").Value = '  +7.23%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '12.06'
$ws.Range("E36").Value = '  +16.08%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.265'
$ws.Range("E37").Value = '  +10.65%  '
$ws.Range("E38").Value = '  +4.14%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '8.620'
$ws.Range("E39").Value = '  +15.59%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.06094'
$ws.Range("E40").Value = '  +4.85%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.02229'
$ws.Range("E41").Value = '  +8.83%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.2021'
$ws.Range("E42").Value = '  +7.53%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.9915'
$ws.Range("E43").Value = '  +3.27%  '
$ws.Range("E44").Value = '  +11.46%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.819'
$ws.Range("E45").Value = '  +7.77%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '13.25'
$ws.Range("E46").Value = '  +7.73%  '
$ws.Range("E47").Value = '  +9.72%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '127.40'
$ws.Range("E48").Value = '  +7.68%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.983'
$ws.Range("E49").Value = '  +8.65%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06814'
$ws.Range("E50").Value = '  +4.95%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '73.90'
$ws.Range("E51").Value = '  +9.22%  '
